$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the tornado-diagram table (Variable / Low / High)
$ws.Range("A11").Value = "Variable"
$ws.Range("B11").Value = "Low"
$ws.Range("C11").Value = "High"

# First column of the data rows
$ws.Range("A12").Value = "Electricity price"
$ws.Range("A13").Value = "Synfuel price"
$ws.Range("A14").Value = "Synfuel plant capacity"

# Remaining header labels (Std low / Std high)
$ws.Range("D11").Value = "Std low"
$ws.Range("E11").Value = "Std high"

# Low/High formula values for each row
$ws.Range("B12").Formula = "=E4-E3"
$ws.Range("C12").Formula = "=E5-E3"
$ws.Range("B13").Formula = "=E6-E3"
$ws.Range("C13").Formula = "=E7-E3"
$ws.Range("B14").Formula = "=E8-E3"
$ws.Range("C14").Formula = "=-E9-E3"

# Column widths for the new Low/High columns (closest achievable to author's bestFit 11.83203125)
$ws.Columns("B:C").ColumnWidth = 11

# Update selection to match the author's final cursor position
$ws.Range("C15").Select()
